# Update cryptos list data (prices / 1h volume changes), and reorder
# three rows (Kaspa / InjectiveProtocol / CoreDAO -> CoreDAO / Kaspa / InjectiveProtocol)
# to match the refreshed data pulled on Mon Apr  1 21:23:37 UTC 2024.
#
# NOTE: every "Price" (column D) value is stored as literal text in the
# workbook (e.g. "1.00", "18.20", "580.82"), even though it looks numeric.
# Assigning a plain numeric-looking string to a Range.Value lets Excel's
# automatic type detection turn it into a real number (dropping trailing
# zeros, introducing floating point noise, etc.), so every Price value is
# written with a leading apostrophe to force it to stay literal text,
# exactly like the original cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 43-45: content swap (new leader CoreDAO moves to row 43) ---
$ws.Range("B43").Value = "CoreDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D43").Value = "'3.53"
$ws.Range("E43").Value = "  +41.42%  "

$ws.Range("B44").Value = "Kaspa"
$ws.Range("C44").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D44").Value = "'0.133"
$ws.Range("E44").Value = "  -2.85%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'34.35"
$ws.Range("E45").Value = "  -6.51%  "

# --- Price (D) / Volume(1h) (E) refresh for every other row ---

$ws.Range("D2").Value = "'69.755.35"
$ws.Range("E2").Value = "  -1.51%  "

$ws.Range("D3").Value = "'3.499.14"
$ws.Range("E3").Value = "  -3.67%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'580.82"
$ws.Range("E5").Value = "  -4.37%  "

$ws.Range("D6").Value = "'193.23"
$ws.Range("E6").Value = "  -2.83%  "

$ws.Range("E7").Value = "  -2.15%  "

$ws.Range("D8").Value = "'3.488.86"
$ws.Range("E8").Value = "  -3.56%  "

$ws.Range("D9").Value = "'1.00"
$ws.Range("E9").Value = "  +0.04%  "

$ws.Range("D10").Value = "'0.204"
$ws.Range("E10").Value = "  -7.81%  "

$ws.Range("D11").Value = "'0.619"
$ws.Range("E11").Value = "  -4.36%  "

$ws.Range("D12").Value = "'51.66"
$ws.Range("E12").Value = "  -4.23%  "

$ws.Range("E13").Value = "  -6.15%  "

$ws.Range("E14").Value = "  -4.31%  "

$ws.Range("D15").Value = "'4.049.39"
$ws.Range("E15").Value = "  -3.85%  "

$ws.Range("D16").Value = "'647.85"
$ws.Range("E16").Value = "  -5.48%  "

$ws.Range("D17").Value = "'69.617.47"
$ws.Range("E17").Value = "  -1.80%  "

$ws.Range("D18").Value = "'3.490.95"
$ws.Range("E18").Value = "  -4.62%  "

$ws.Range("E19").Value = "  -4.76%  "

$ws.Range("E20").Value = "  -1.70%  "

$ws.Range("D21").Value = "'18.25"
$ws.Range("E21").Value = "  -4.05%  "

$ws.Range("E22").Value = "  -5.05%  "

$ws.Range("D23").Value = "'18.20"
$ws.Range("E23").Value = "  -3.29%  "

$ws.Range("D24").Value = "'5.24"
$ws.Range("E24").Value = "  -3.31%  "

$ws.Range("D25").Value = "'98.78"
$ws.Range("E25").Value = "  -5.96%  "

$ws.Range("E26").Value = "  -7.54%  "

$ws.Range("E27").Value = "  -3.82%  "

$ws.Range("E28").Value = "  -3.90%  "

$ws.Range("D29").Value = "'9.34"
$ws.Range("E29").Value = "  -5.36%  "

$ws.Range("D30").Value = "'32.65"
$ws.Range("E30").Value = "  -4.89%  "

$ws.Range("D31").Value = "'4.24"
$ws.Range("E31").Value = "  -7.87%  "

$ws.Range("D32").Value = "'6.74"
$ws.Range("E32").Value = "  -6.04%  "

$ws.Range("E33").Value = "  -4.39%  "

$ws.Range("E34").Value = "  -4.81%  "

$ws.Range("D35").Value = "'61.34"
$ws.Range("E35").Value = "  -3.12%  "

$ws.Range("D36").Value = "'527.99"
$ws.Range("E36").Value = "  +4.82%  "

$ws.Range("D37").Value = "'3.707.80"
$ws.Range("E37").Value = "  -6.05%  "

$ws.Range("E38").Value = "  +0.30%  "

$ws.Range("E39").Value = "  -9.38%  "

$ws.Range("D40").Value = "'3.55"
$ws.Range("E40").Value = "  -0.11%  "

$ws.Range("D41").Value = "'2.91"
$ws.Range("E41").Value = "  -4.28%  "

$ws.Range("E42").Value = "  -3.76%  "

$ws.Range("E46").Value = "  -3.23%  "

$ws.Range("D47").Value = "'3.38"
$ws.Range("E47").Value = "  -3.05%  "

$ws.Range("E48").Value = "  -7.70%  "

$ws.Range("E49").Value = "  -4.17%  "

$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "  -0.30%  "

$ws.Range("E51").Value = "  -5.49%  "
